# Vendors.xlsx — restock count updates after a shopping run, plus
# re-selecting the "Grocery" tab as the active sheet (was "Bath").
#
# Commit message: "Can write a ShopsToVisit.txt" — the stock levels below
# reflect items that were bought (decremented from the default starting
# stock of 5) so a downstream script can flag which shops still need a
# visit.

$wb = $excel.ActiveWorkbook

# --- Grocery sheet: Stock (column B) updates ---
$groceryWs = $wb.Worksheets.Item("Grocery")
$groceryWs.Range("B2").Value = 4   # Pizza
$groceryWs.Range("B3").Value = 3   # Fruits
$groceryWs.Range("B4").Value = 3   # Vegetables
$groceryWs.Range("B6").Value = 4   # Basic Dog Food

# --- Pet sheet: Stock (column B) updates ---
$petWs = $wb.Worksheets.Item("Pet")
$petWs.Range("B2").Value = 4   # Brush
$petWs.Range("B3").Value = 3   # Premium Cat Food
$petWs.Range("B5").Value = 4   # Cat Litter

# --- Bath sheet: Stock (column B) updates ---
$bathWs = $wb.Worksheets.Item("Bath")
$bathWs.Range("B2").Value = 4   # Body Butter
$bathWs.Range("B3").Value = 4   # Oatmeal Soap

# --- Make "Grocery" the active/selected tab (was "Bath") ---
$groceryWs.Activate()
